$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

# The sheet currently has two mailto: hyperlinks (B3/B4) rendered with the
# built-in "Hyperlink" cell style. The new data no longer needs them, so
# drop the hyperlinks and restore those cells to the default style before
# wiping the old table.
$ws.Cells.Hyperlinks.Delete()
$ws.Range("B3:B4").Style = "Normal"

# Clear out the whole previous user table (A2:C15 -> up to 15 rows incl. a
# "Result" column) so we can lay down the fresh data beneath the header row.
$ws.Range("A2:C15").Clear()

# Remove the now-unused "Hyperlink" named cell style from the workbook.
$wb.Styles.Item("Hyperlink").Delete()

# New user data: just Name / E-mail pairs, checking what's already present.
$data = @(
    @("Olga",  "olga@gmail.com"),
    @("Una",   "una@gmail.com"),
    @("Tanja", "tanja@gmail.com"),
    @("Olga",  "olga@gmail.com"),
    @("Olga",  "olga@gmail.com"),
    @("Lena",  "lena@gmail.com"),
    @("Olga",  "olga@gmail.com"),
    @("Olga",  "olga@gmail.com"),
    @("Lena",  "lena@gmail.com"),
    @("Tanja", "tanja@gmail.com")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

$ws.Range("D10").Select()

$wb.Save()
